$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H41").Value = 509.28
$ws.Range("I41").Value = 342
$ws.Range("J41").Value = 620.8
$ws.Range("K41").Value = 342
$ws.Range("L41").Value = 620.8
$ws.Range("M41").Value = 98
$ws.Range("N41").Value = -1500.8
$ws.Range("H44").Value = 11344.444
$ws.Range("I44").Value = 100
$ws.Range("J44").Value = 12750
$ws.Range("K44").Value = 100
$ws.Range("L44").Value = 12750
$ws.Range("M44").Value = 362
$ws.Range("N44").Value = -13674
$ws.Range("H58").Value = 3160
$ws.Range("I58").Value = 400
$ws.Range("J58").Value = 5000
$ws.Range("K58").Value = 1200
$ws.Range("L58").Value = 15000
$ws.Range("M58").Value = -1050
$ws.Range("N58").Value = -15300
$ws.Range("H74").Value = 4500
$ws.Range("I74").Value = 4000
$ws.Range("K74").Value = 4000
$ws.Range("M74").Value = -3064
$ws.Range("H77").Value = 4500
$ws.Range("I77").Value = 4000
$ws.Range("K77").Value = 20000
$ws.Range("M77").Value = -15320
$ws.Range("H116").Value = 8351777
$ws.Range("I116").Value = 13916207
$ws.Range("K116").Value = 13916207
$ws.Range("M116").Value = -13912765
$ws.Range("H138").Value = 247065.4
$ws.Range("I138").Value = 425831.56
$ws.Range("J138").Value = 6119.7393
$ws.Range("K138").Value = 1277494.68
$ws.Range("L138").Value = 18359.2179
$ws.Range("M138").Value = -1272354.68
$ws.Range("N138").Value = -28639.2179

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H18").Value = 27999
$ws.Range("J18").Value = 27999
$ws.Range("L18").Value = 27999
$ws.Range("N18").Value = -28643
$ws.Range("H31").Value = 4209.8
$ws.Range("I31").Value = 4209.8
$ws.Range("K31").Value = 4209.8
$ws.Range("M31").Value = -3915.8
$ws.Range("H32").Value = 7437.0713
$ws.Range("I32").Value = 7437.0713
$ws.Range("K32").Value = 7437.0713
$ws.Range("M32").Value = -7150.0713
$ws.Range("H45").Value = 6314.7
$ws.Range("I45").Value = 5580.875
$ws.Range("K45").Value = 5580.875
$ws.Range("M45").Value = -5203.875
$ws.Range("H74").Value = 3077.6582
$ws.Range("I74").Value = 6028.0835
$ws.Range("J74").Value = 1790.2
$ws.Range("K74").Value = 6028.0835
$ws.Range("L74").Value = 1790.2
$ws.Range("M74").Value = -5154.0835
$ws.Range("N74").Value = -3538.2
$ws.Range("H77").Value = 3077.6582
$ws.Range("I77").Value = 6028.0835
$ws.Range("J77").Value = 1790.2
$ws.Range("K77").Value = 30140.4175
$ws.Range("L77").Value = 8951
$ws.Range("M77").Value = -25772.4175
$ws.Range("N77").Value = -17687
$ws.Range("H97").Value = 66701668
$ws.Range("J97").Value = 100002500
$ws.Range("L97").Value = 100002500
$ws.Range("N97").Value = -100003492
$ws.Range("H139").Value = 166578.5
$ws.Range("J139").Value = 166578.5
$ws.Range("L139").Value = 166578.5
$ws.Range("N139").Value = -176858.5

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H5").Value = 2948.5
$ws.Range("I5").Value = 1649.5
$ws.Range("J5").Value = 4247.5
$ws.Range("K5").Value = 1649.5
$ws.Range("L5").Value = 4247.5
$ws.Range("M5").Value = -1536.5
$ws.Range("N5").Value = -4473.5
$ws.Range("H7").Value = 8169573.5
$ws.Range("I7").Value = 19187.375
$ws.Range("J7").Value = 14097127
$ws.Range("K7").Value = 19187.375
$ws.Range("L7").Value = 14097127
$ws.Range("M7").Value = -19074.375
$ws.Range("N7").Value = -14097353
$ws.Range("H44").Value = 24500
$ws.Range("I44").Value = 0
$ws.Range("J44").Value = 24500
$ws.Range("K44").Value = 0
$ws.Range("L44").Value = 24500
$ws.Range("M44").ClearContents()
$ws.Range("N44").Value = -25494
$ws.Range("H134").Value = 2007.081
$ws.Range("I134").Value = 1435.2122
$ws.Range("J134").Value = 6725
$ws.Range("K134").Value = 4305.6366
$ws.Range("L134").Value = 20175
$ws.Range("M134").Value = -1770.6366
$ws.Range("N134").Value = -25245

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H12").Value = 14995.5
$ws.Range("J12").Value = 14995.5
$ws.Range("L12").Value = 14995.5
$ws.Range("N12").Value = -15335.5
$ws.Range("H15").Value = 9138.727999999999
$ws.Range("H31").Value = 3292.946
$ws.Range("I31").Value = 2196
$ws.Range("J31").Value = 4583.4707
$ws.Range("K31").Value = 2196
$ws.Range("L31").Value = 4583.4707
$ws.Range("M31").Value = -1901
$ws.Range("N31").Value = -5173.4707
$ws.Range("H34").Value = 3292.946
$ws.Range("I34").Value = 2196
$ws.Range("J34").Value = 4583.4707
$ws.Range("K34").Value = 2196
$ws.Range("L34").Value = 4583.4707
$ws.Range("M34").Value = -1994
$ws.Range("N34").Value = -4987.4707
$ws.Range("H35").Value = 6161.25
$ws.Range("I35").Value = 6733.6665
$ws.Range("K35").Value = 6733.6665
$ws.Range("M35").Value = -6439.6665
$ws.Range("H105").Value = 3938.2058
$ws.Range("I105").Value = 7049.222
$ws.Range("J105").Value = 438.3125
$ws.Range("K105").Value = 7049.222
$ws.Range("L105").Value = 438.3125
$ws.Range("M105").Value = -5302.222
$ws.Range("N105").Value = -3932.3125
$ws.Range("H141").Value = 313707.94
$ws.Range("J141").Value = 349920.3
$ws.Range("L141").Value = 349920.3
$ws.Range("N141").Value = -360280.3

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 627520.5
$ws.Range("I5").Value = 1041.75
$ws.Range("J5").Value = 1253999.2
$ws.Range("K5").Value = 3125.25
$ws.Range("L5").Value = 3761997.6
$ws.Range("M5").Value = -3013.25
$ws.Range("N5").Value = -3762221.6
$ws.Range("H7").Value = 123.09091
$ws.Range("J7").Value = 157.25
$ws.Range("L7").Value = 471.75
$ws.Range("N7").Value = -695.75
$ws.Range("H107").Value = 1467.2632
$ws.Range("I107").Value = 441.66666
$ws.Range("K107").Value = 1324.99998
$ws.Range("M107").Value = 595.0000199999999
$ws.Range("H113").Value = 1210.1666
$ws.Range("I113").Value = 916.5
$ws.Range("J113").Value = 1268.9
$ws.Range("K113").Value = 2749.5
$ws.Range("L113").Value = 3806.7
$ws.Range("M113").Value = -579.5
$ws.Range("N113").Value = -8146.700000000001
$ws.Range("H121").Value = 1348.5883
$ws.Range("J121").Value = 1843.6
$ws.Range("L121").Value = 5530.799999999999
$ws.Range("N121").Value = -8150.799999999999
$ws.Range("H131").Value = 23812684
$ws.Range("I131").Value = 45458796
$ws.Range("J131").Value = 1960.05
$ws.Range("K131").Value = 136376388
$ws.Range("L131").Value = 5880.15
$ws.Range("M131").Value = -136371348
$ws.Range("N131").Value = -15960.15
$ws.Range("H135").Value = 627520.5
$ws.Range("I135").Value = 1041.75
$ws.Range("J135").Value = 1253999.2
$ws.Range("K135").Value = 9375.75
$ws.Range("L135").Value = 11285992.8
$ws.Range("M135").Value = -6840.75
$ws.Range("N135").Value = -11291062.8

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 7152.0435
$ws.Range("I16").Value = 7309.381
$ws.Range("K16").Value = 7309.381
$ws.Range("M16").Value = -7139.381
$ws.Range("H20").Value = 15000
$ws.Range("I20").Value = 15000
$ws.Range("K20").Value = 15000
$ws.Range("M20").Value = -14774
$ws.Range("H68").Value = 5216.905
$ws.Range("J68").Value = 6361.143
$ws.Range("L68").Value = 6361.143
$ws.Range("N68").Value = -7859.143
$ws.Range("H71").Value = 5216.905
$ws.Range("J71").Value = 6361.143
$ws.Range("L71").Value = 31805.715
$ws.Range("N71").Value = -39293.715
$ws.Range("H82").Value = 1845.4667
$ws.Range("J82").Value = 1914
$ws.Range("L82").Value = 1914
$ws.Range("N82").Value = -2636
$ws.Range("H85").Value = 1845.4667
$ws.Range("J85").Value = 1914
$ws.Range("L85").Value = 1914
$ws.Range("N85").Value = -4410
$ws.Range("H104").Value = 30369
$ws.Range("J104").Value = 30369
$ws.Range("L104").Value = 30369
$ws.Range("N104").Value = -37357

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 443613.53
$ws.Range("I136").Value = 516825.8
$ws.Range("J136").Value = 4340
$ws.Range("K136").Value = 1550477.4
$ws.Range("L136").Value = 13020
$ws.Range("M136").Value = -1547927.4
$ws.Range("N136").Value = -18120
